$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'51.661.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.33%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.890.10"
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'352.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.07%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'108.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.45%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.559"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.59%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.20%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.15%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'38.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.57%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "'  +1.01%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.0861"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.60%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'19.29"
$ws.Range("D13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'3.350.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.74%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'7.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.31%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'2.869.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.29%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'0.965"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -4.21%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'51.593.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.50%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'3.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.34%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'7.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.63%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'13.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.96%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.0₃0967"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.50%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'69.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.77%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'266.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.80%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'2.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.18%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'  +8.56%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'26.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.88%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("B28").Value = "'Dai"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.05%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("B29").Value = "'Filecoin"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'7.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +13.41%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  +7.80%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'10.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.55%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'37.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.93%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("E33").Value = "'  -2.53%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("E34").Value = "'  -3.62%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'51.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.36%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'0.0434"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.75%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'  +0.16%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'3.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.50%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'17.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -4.24%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'1.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.02%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "'  -6.37%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'  +0.53%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'22.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -5.36%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'118.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.32%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'2.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.83%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'2.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -5.56%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("B47").Value = "'Maker"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'2.112.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.98%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("B48").Value = "'NEARProtocol"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'3.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.47%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  -7.65%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.0334"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.51%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("B51").Value = "'MultiversX"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'61.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.75%  "
$ws.Range("E51").Style = "Normal"

